# Weekly update: a new daily price record was added for
# "Terminal Hortofrutícola Agro Chillán - Zapallo italiano".
# In the source table this new record is inserted right before the
# existing row 171, pushing all subsequent rows (171..245) down by one
# (to 172..246), and the sheet's used range grows from A1:R245 to A1:R246.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 171; everything from 171 downward shifts
# down by one row (old row 171 becomes 172, ..., old row 245 becomes 246).
$ws.Rows("171:171").Insert()

# Populate the newly inserted row 171 with the new weekly record.
$ws.Range("A171").Value = 7
$ws.Range("B171").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C171").Value = "Ñuble"
$ws.Range("D171").Value = 44825
$ws.Range("E171").Value = 16
$ws.Range("F171").Value = 100112032
$ws.Range("G171").Value = "Zapallo italiano"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 120
$ws.Range("K171").Value = 14000
$ws.Range("L171").Value = 15000
$ws.Range("M171").Value = 14500
$ws.Range("N171").Value = "`$/caja 50 unidades"
$ws.Range("O171").Value = "Región de Arica y Parinacota"
$ws.Range("P171").Value = 290
$ws.Range("Q171").Value = 50
$ws.Range("R171").Value = "Hortaliza"
